$d = $word.ActiveDocument

function Set-ParagraphOoxml($AnchorText, $InnerXml) {
    $full = $d.Content
    $full.Find.ClearFormatting()
    $full.Find.Execute($AnchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $p = $full.Paragraphs(1)
    $pr = $p.Range
    $start = $pr.Start
    $end = $pr.End
    $textRange = $d.Range($start, $end)
    $ooxml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
$InnerXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $textRange.InsertXML($ooxml) | Out-Null
}

# 1) "Wear_Rate" gets wrapped in spellStart/spellEnd proofErr markers, splitting the
#    surrounding run into three runs.
Set-ParagraphOoxml "Our journey began by importing the training dataset" @'
          <w:p w14:paraId="05119501" w14:textId="77777777" w:rsidR="00F00999" w:rsidRDefault="00F00999" w:rsidP="00F00999">
            <w:r><w:t xml:space="preserve">Our journey began by importing the training dataset, &quot;training_data.csv.&quot; This dataset consisted of RPM, Load, Hardness, and </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Wear_Rate</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> columns. Fortunately, the data was relatively clean and did not require extensive preprocessing, which allowed us to focus on the modeling process.</w:t></w:r>
          </w:p>
'@

# 2) The markdown image syntax "![" gets flagged with a gramStart/gramEnd proofErr pair
#    (gramStart placed before the run, per the diff).
Set-ParagraphOoxml "Polynomial Order vs. Average RMSE" @'
          <w:p w14:paraId="43F262D4" w14:textId="18190824" w:rsidR="00F00999" w:rsidRDefault="00F00999" w:rsidP="00F00999">
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>![</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t>Polynomial Order vs. Average RMSE](polynomial_order_vs_rmse.png)</w:t></w:r>
          </w:p>
'@

# 3) "one fold" gets wrapped in gramStart/gramEnd proofErr markers.
Set-ParagraphOoxml "The principle of 5-fold cross-validation" @'
          <w:p w14:paraId="77A17971" w14:textId="77777777" w:rsidR="00F00999" w:rsidRDefault="00F00999" w:rsidP="00F00999">
            <w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">The principle of 5-fold cross-validation entails splitting the data into five subsets, known as folds. During each iteration, </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>one fold</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> is designated as the validation set, while the remaining four folds serve as the training data. RMSE is calculated for each fold, and the average RMSE across all folds is reported. This technique ensures that our models' predictive performance is evaluated comprehensively and helps us assess their generalizability to unseen data.</w:t></w:r>
          </w:p>
'@

# 4) "joblib" gets wrapped in spellStart/spellEnd proofErr markers.
Set-ParagraphOoxml "we used Python's joblib library" @'
          <w:p w14:paraId="2F3673FE" w14:textId="77777777" w:rsidR="00F00999" w:rsidRDefault="00F00999" w:rsidP="00F00999">
            <w:r><w:t xml:space="preserve">To ensure that our models can be evaluated on a separate test dataset, we used Python&apos;s </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>joblib</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> library to save the best-performing models from both Part 1 and Part 2. These saved models can undergo rigorous testing to validate their real-world predictive performance.</w:t></w:r>
          </w:p>
'@

# 5) Append " CHANGE" (red) after the Part 2 "### Results" heading.
$full = $d.Content
$full.Find.ClearFormatting()
$full.Find.Execute("### Polynomial Models with Ridge Regularization", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$after = $d.Range($full.End, $d.Content.End)
$after.Find.ClearFormatting()
$after.Find.Execute("### Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$after.Collapse(0)
$after.InsertAfter(" ")
$after.Collapse(0)
$changeRange = $d.Range($after.End, $after.End)
$changeRange.InsertAfter("CHANGE")
$changeRange.Font.Color = 255
